$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$x = $r.XML
Write-Output $x.GetType()
